$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# C7 was stored as a text value "5" — convert it to a true number.
$ws.Range("C7").Value = 5

# New submission row synced in.
$ws.Range("A8").Value = "2026-02-08 20:20:44"
$ws.Range("B8").Value = "Rachel Michael"

# C8 ("Admission No") stays text, same as how some prior rows store it —
# force text formatting before assigning a numeric-looking string, then
# drop the formatting again so no extra style is left on the cell.
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "34"
$ws.Range("C8").ClearFormats()

$ws.Range("D8").Value = 10
